$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Columns B (ORDER_DATE) and C (DELIVERY_DATE): switch from a date
# number format to plain text, then rewrite the values as short-form text
# dates (matches the author's "make sure TEXT format (MM/DD/RRRR)" note). ---
$ws.Range("B1:C4").NumberFormat = "@"

$ws.Range("B2").Value = "9/13/25"
$ws.Range("B3").Value = "9/13/25"
$ws.Range("B4").Value = "9/13/25"

$ws.Range("C2").Value = "9/20/25"
$ws.Range("C3").Value = "9/20/25"
$ws.Range("C4").Value = "9/20/25"

# --- Column D (ORDER_STATUS): 3 -> 2 (still a valid item in the existing
# "2,3,5" data-validation list). ---
$ws.Range("D2").Value = "2"
$ws.Range("D3").Value = "2"
$ws.Range("D4").Value = "2"

# --- Column G (QTY): 5 -> 45. ---
$ws.Range("G2").Value = "45"
$ws.Range("G3").Value = "45"
$ws.Range("G4").Value = "45"

# --- Comments: add new threaded comments on B1/C1 reminding to keep the
# TEXT format, and drop the stale "Optional." threaded comment that used
# to live on F1, replacing it with a blank legacy comment. ---
$excel.UserName = "Author"

$ws.Range("B1").AddCommentThreaded("make sure`n TEXT format (MM/DD/RRRR)")
$ws.Range("C1").AddCommentThreaded("make sure`n TEXT format (MM/DD/RRRR)")

$ws.Range("F1").CommentThreaded.Delete()
$ws.Range("F1").AddComment()

# --- Selection cosmetic change. ---
$ws.Range("F18").Select()
